$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseSteps")

# Clear the existing data rows (rows 2-8) and rebuild the table with the
# new, smaller test-case set (login module + Admin user management).
$ws.Rows("2:8").ClearContents()

$ws.Range("A1").Value = "TestCaseID"
$ws.Range("B1").Value = "ClassName"
$ws.Range("C1").Value = "steps"
$ws.Range("D1").Value = "ActiveFlag"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Login_Page"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Add_User"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1

$ws.Range("D9").Select()
